$wb = $excel.ActiveWorkbook

# --- library_content sheet (sheet1.xml) ---
$ws1 = $wb.Worksheets.Item("library_content")

# fix NIST CSF score 1: bump library_version from 1 to 2
$ws1.Cells.Item(2, 2).Value2 = 2

# insert two new rows (14 & 15) before the existing "tab" rows to hold the
# new framework_min_score / framework_max_score entries; everything below
# shifts down automatically (old row 14 -> 16, old row 15 -> 17)
$ws1.Rows.Item(14).Insert()
$ws1.Rows.Item(14).Insert()

$ws1.Cells.Item(14, 1).Value2 = "framework_min_score"
$ws1.Cells.Item(14, 2).Value2 = 1

$ws1.Cells.Item(15, 1).Value2 = "framework_max_score"
$ws1.Cells.Item(15, 2).Value2 = 4

# --- scores sheet (sheet3.xml) ---
$ws3 = $wb.Worksheets.Item("scores")
$ws3.Range("B4").Select()

# re-select the library_content sheet so it remains the active tab,
# matching the updated selection there
$ws1.Range("B19").Select()
